$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header C1 from "Correct" to "Correct Answer"
$ws.Range("C1").Value = "Correct Answer"

# Add new question row (row 3)
$ws.Range("A3").Value = $true
$ws.Range("B3").Value = "What is the meaning PAAU"
$ws.Range("C3").Value = "b"
$ws.Range("D3").Value = "a: Prince Adeiza Audu University;b: Prince Abubakar Audu University;c: Prince Abdullahi Audu University;d: Prince Adeiza Ahmodu University;"
$ws.Range("E3").Value = 10

# Widen columns C and D to fit new content
$ws.Columns.Item(3).ColumnWidth = 14.583333333333332
$ws.Columns.Item(4).ColumnWidth = 124.58333333333334

# Update the view: scroll to column D and select D3
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("D3").Select() | Out-Null
